$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Apply the "Table Grid" style to the table (adds <w:tblStyle w:val="TableGrid"/>)
$t.Style = "Table Grid"

# Update the measured results in the "t" / "t_max" columns (columns 12 & 13)
$t.Cell(2, 12).Range.Text  = "0.0038"
$t.Cell(2, 13).Range.Text  = "0.0206"
$t.Cell(3, 13).Range.Text  = "0.0076"
$t.Cell(4, 12).Range.Text  = "0.0036"
$t.Cell(4, 13).Range.Text  = "0.0178"
$t.Cell(5, 13).Range.Text  = "0.0161"
$t.Cell(6, 13).Range.Text  = "0.0109"
$t.Cell(7, 12).Range.Text  = "0.0038"
$t.Cell(7, 13).Range.Text  = "0.0099"

# Append a new, empty trailing row (13 blank cells) to the table
$t.Rows.Add() | Out-Null
